$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.02, 1.039877732955504, 1.046959125648108, 1.038345866015583, 1.057589790228187, 1.042073981240425, 1.044967369470309, 1.049723102938305, 1.041134161943574, 1.060324344721522, 1.018872535769837),
    @(1.02, 1.040817298629185, 1.047690009615101, 1.039145120743924, 1.058483129053316, 1.042319073527953, 1.045552371696296, 1.050265944494896, 1.041743439435209, 1.061031347811438, 1.019070088669788),
    @(1.02, 1.041425774192972, 1.048163405259247, 1.039663089298199, 1.059062024736793, 1.042476774209586, 1.045930793120266, 1.050616994179363, 1.042137825578336, 1.061489033655254, 1.019197797657611),
    @(1.02, 1.041681699233354, 1.048362530650741, 1.039881032725264, 1.05930559343237, 1.04254285782449, 1.046089853140616, 1.050764525320275, 1.042303658613015, 1.061681493197563, 1.019251457139004),
    @(1.02, 1.041724677291583, 1.048395971081235, 1.039917637453408, 1.05934650142782, 1.042553941002789, 1.046116558326014, 1.050789293477053, 1.042331504618234, 1.061713810803438, 1.019260465067098),
    @(1.02, 1.041429193398189, 1.048166065552572, 1.039666000726332, 1.059065278525928, 1.042477658062706, 1.045932918600701, 1.050618965696037, 1.04214004131878, 1.061491605120447, 1.019198514774188),
    @(1.02, 1.040195156549579, 1.047206034202516, 1.038615811947791, 1.057891522239609, 1.042156995341517, 1.045165096728609, 1.049906600622972, 1.041340039880495, 1.0605632359781, 1.018939324486875),
    @(1.02, 1.03802460834186, 1.045517959285457, 1.036771414627884, 1.055829758102093, 1.041585152984066, 1.043811267542923, 1.048649797694991, 1.039931478003686, 1.058928977197268, 1.018481690388922),
    @(1.02, 1.036580316999344, 1.044395096196326, 1.035546049077581, 1.054459738764616, 1.041199391669215, 1.042908217813033, 1.047810965812324, 1.038993270349543, 1.057840653935637, 1.018176015652179),
    @(1.02, 1.035955585972228, 1.043909499487477, 1.035016473908872, 1.053867589068432, 1.041031284772311, 1.042517082001855, 1.047447526763345, 1.038587227392703, 1.057369694255337, 1.018043521487122),
    @(1.02, 1.035723632708532, 1.043729220365211, 1.0348199201993, 1.053647801719902, 1.04096868222937, 1.042371781307039, 1.047312497384271, 1.038436437156122, 1.05719480390316, 1.017994287385003),
    @(1.02, 1.035773382945207, 1.043767886612576, 1.03486207464172, 1.053694939433325, 1.040982117925832, 1.042402949482135, 1.047341463073231, 1.038468780699499, 1.057232316418989, 1.018004849159603),
    @(1.02, 1.035936410604331, 1.043894595651297, 1.035000223566403, 1.053849418035364, 1.041026113290511, 1.042505071710346, 1.047436365838321, 1.038574762361685, 1.057355236838166, 1.018039452183849),
    @(1.02, 1.036036870561131, 1.043972677652955, 1.035085362138795, 1.053944619099348, 1.041053199081057, 1.042567990569254, 1.047494834362868, 1.038640065460139, 1.057430978171366, 1.018060769628451),
    @(1.02, 1.03662179226038, 1.044427336631624, 1.035581216768345, 1.054499060623954, 1.041210525880164, 1.042934173989424, 1.047835081531022, 1.039020222539424, 1.057871916226385, 1.018184806044559),
    @(1.02, 1.036988874698902, 1.044712696617371, 1.035892526341533, 1.0548471367376, 1.04130892682549, 1.043163842601151, 1.048048451590916, 1.039258741249633, 1.058148583937862, 1.018262574961201),
    @(1.02, 1.037203050965007, 1.044879201071242, 1.03607420590237, 1.055050267529201, 1.041366219184048, 1.043297793802841, 1.048172885508715, 1.039397885092564, 1.058309987572053, 1.01830792320273),
    @(1.02, 1.037276090241696, 1.044935984718169, 1.036136170509212, 1.055119547454624, 1.041385736868272, 1.043343465870789, 1.048215310627514, 1.039445332892394, 1.058365026758098, 1.018323383554893),
    @(1.02, 1.036949483673618, 1.044682074099166, 1.035859115628493, 1.054809780708771, 1.041298380010737, 1.043139202423391, 1.0480255611912, 1.039233148381703, 1.058118897201528, 1.018254232439484),
    @(1.02, 1.035888400263018, 1.04385728041575, 1.034959537886424, 1.053803923414264, 1.041013162168807, 1.042474999677747, 1.047408420212775, 1.038543552506271, 1.05731903858497, 1.018029263000877),
    @(1.02, 1.035221832086793, 1.043339239332155, 1.034394829858895, 1.05317244721795, 1.040832908051701, 1.042057300211972, 1.047020214927698, 1.038110163022429, 1.056816396950065, 1.017887701027853),
    @(1.02, 1.035575137452758, 1.043613811109188, 1.034694107225827, 1.053507114506938, 1.040928551790229, 1.042278738693276, 1.047226027020691, 1.038339892920807, 1.057082831611474, 1.017962756448867),
    @(1.02, 1.036967282586895, 1.044695910915481, 1.035874212189065, 1.054826659969408, 1.041303145981413, 1.043150336291988, 1.048035904443221, 1.039244712635667, 1.058132311273323, 1.01825800210632),
    @(1.02, 1.038585268649189, 1.045953929086853, 1.037247495708751, 1.056361988739178, 1.041733789469476, 1.044161356857018, 1.048974885691925, 1.040295482862591, 1.0593512693773, 1.018600104769281)
)

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)
$startRow = 2

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $cols[$j]).Value = $row[$j]
    }
}
